$wb = $excel.ActiveWorkbook

# Rename the formula "IsLastOfMonth" -> "isLastOfMonth" on the Library_Formula sheet
$wsFormula = $wb.Worksheets.Item("Library_Formula")
$wsFormula.Range("C2").Value = "isLastOfMonth"

# Update the selection remembered on the (now inactive) "Library" sheet
$wsLibrary = $wb.Worksheets.Item("Library")
$wsLibrary.Activate()
$wsLibrary.Range("C5").Select()

# Make "Library_Formula" the active tab with its own updated selection
$wsFormula.Activate()
$wsFormula.Range("C3").Select()
